$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.080.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.602.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.05%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.486"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.247"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0613"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("E11").Value = "  +4.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.602.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.074.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.70%  "
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("E24").Value = "  +12.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.122"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0472"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.125.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.70%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0164"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.34%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.785"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.51%  "
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.783"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.738.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0504"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₇0921"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -16.86%  "
